# Symbol-list refresh (GitHub Actions update, Fri Jan 27 05:23:05 UTC 2023).
# Updates scraped Price/Volume(1h) figures for the existing coin rows and
# corrects the coin identity + link that had landed on the wrong two rows
# (row 20 <-> row 21: ProBitToken vs MCDex).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / Link swap (rows 20-21) ---
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

# --- Price (D) / Volume(1h) (E) refresh ---
# Price and Volume(1h) are stored as literal TEXT in this sheet (e.g. "303.00",
# "-1.40%"), not numbers. Excel normally auto-converts a plain numeric-looking
# string assignment into a real Number/percentage, so each value below is
# entered with a leading apostrophe (the standard "force text" input prefix)
# to keep it a literal text cell, matching the scraped source data.

$ws.Range("D2").Value = "'302.64"
$ws.Range("E2").Value = "'-1.54%"
$ws.Range("E3").Value = "'-1.61%"
$ws.Range("D4").Value = "'5.033"
$ws.Range("D5").Value = "'0.07896"
$ws.Range("E5").Value = "'-2.84%"
$ws.Range("E6").Value = "'-5.13%"
$ws.Range("D7").Value = "'7.785"
$ws.Range("E7").Value = "'0.06%"
$ws.Range("D8").Value = "'0.9194"
$ws.Range("E8").Value = "'-1.23%"
$ws.Range("D9").Value = "'0.1342"
$ws.Range("E9").Value = "'-4.00%"
$ws.Range("D10").Value = "'0.1887"
$ws.Range("D11").Value = "'0.09064"
$ws.Range("E11").Value = "'-1.91%"
$ws.Range("D12").Value = "'0.03467"
$ws.Range("E12").Value = "'1.29%"
$ws.Range("D13").Value = "'0.09799"
$ws.Range("E13").Value = "'-0.73%"
$ws.Range("D14").Value = "'0.001407"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("D15").Value = "'0.006032"
$ws.Range("E15").Value = "'4.77%"
$ws.Range("E16").Value = "'3.40%"
$ws.Range("D17").Value = "'4.105"
$ws.Range("E17").Value = "'-1.92%"
$ws.Range("D18").Value = "'3.303"
$ws.Range("D19").Value = "'0.3437"
$ws.Range("E19").Value = "'-0.04%"
$ws.Range("D20").Value = "'5.157"
$ws.Range("E20").Value = "'5.31%"
$ws.Range("D21").Value = "'0.1301"
$ws.Range("E21").Value = "'-1.61%"
$ws.Range("D22").Value = "'0.2193"
$ws.Range("E22").Value = "'-12.27%"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("E23").Value = "'-2.22%"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("D25").Value = "'0.004601"
$ws.Range("E25").Value = "'-5.47%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'4.88%"
$ws.Range("D27").Value = "'0.0004444"
$ws.Range("E27").Value = "'0.05%"
$ws.Range("D39").Value = "'0.01929"
$ws.Range("E39").Value = "'-3.68%"
$ws.Range("D40").Value = "'0.05255"
$ws.Range("E40").Value = "'6.36%"
$ws.Range("D41").Value = "'0.007606"
$ws.Range("E41").Value = "'-0.73%"
$ws.Range("D42").Value = "'0.01015"
$ws.Range("E42").Value = "'-0.40%"
$ws.Range("E43").Value = "'-2.94%"
$ws.Range("D44").Value = "'0.002162"
$ws.Range("E44").Value = "'2.89%"
$ws.Range("E45").Value = "'-1.62%"
$ws.Range("D46").Value = "'0.00006152"
$ws.Range("E46").Value = "'-4.51%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("D49").Value = "'0.001660"
$ws.Range("E49").Value = "'39.34%"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E51").Value = "'0.00%"
